$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Stomatitis row (row 10) treatment/AE cost data
$ws.Range("C10").Value = 940
$ws.Range("D10").Value = 106
$ws.Range("E10").Value = 1844

# Update the active cell selection to reflect where editing ended
$ws.Range("E10").Select()
